# This edit reshuffles the data rows (2-16) of the "Artfynd" worksheet:
# the content of each row moves to a different row position (a pure
# permutation of whole rows), while rows 1 (header) and 17-19 stay put.
#
# Mapping: new row N gets the content that used to live in row $map[N].
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{}
$map[2]  = 6
$map[3]  = 14
$map[4]  = 15
$map[5]  = 16
$map[6]  = 2
$map[7]  = 3
$map[8]  = 4
$map[9]  = 5
$map[10] = 7
$map[11] = 8
$map[12] = 9
$map[13] = 10
$map[14] = 11
$map[15] = 12
$map[16] = 13

$firstRow = 2
$lastRow = 16
$lastCol = "AY"

# Make sure text-typed columns that merely look numeric/date-like keep
# being stored as plain text instead of being auto-converted to actual
# numbers / Excel date serials when we write the values back:
#  - I: "Antal" holds small integers stored as text
#  - Y/AA: "Startdatum"/"Slutdatum" hold dates stored as text
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Range("I" + $r).NumberFormat = "@"
    $ws.Range("Y" + $r).NumberFormat = "@"
    $ws.Range("AA" + $r).NumberFormat = "@"
}

# Snapshot every source row BEFORE any writes happen, since several rows
# both give and receive content in this permutation.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $snapshot[$r] = $ws.Range("A" + $r + ":" + $lastCol + $r).Value()
}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $srcRow = $map[$r]
    $ws.Range("A" + $r + ":" + $lastCol + $r).Value = $snapshot[$srcRow]
}
